$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.042.81'
$ws.Range("E2").Value = '  -7.09%  '

$ws.Range("D3").Value = '1.417.74'
$ws.Range("E3").Value = '  -7.47%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.36%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9985'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '274.75'
$ws.Range("E6").Value = '  -5.09%  '

$ws.Range("E7").Value = '  -5.54%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3130'
$ws.Range("E8").Value = '  -1.34%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.80'
$ws.Range("E9").Value = '  -7.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.043'
$ws.Range("E10").Value = '  -2.39%  '

$ws.Range("E11").Value = '  -9.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9996'
$ws.Range("E12").Value = '  -0.36%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.493'
$ws.Range("E13").Value = '  -4.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.73'
$ws.Range("E14").Value = '  -2.44%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.202'
$ws.Range("E15").Value = '  -5.76%  '

$ws.Range("D16").Value = '1.418.30'
$ws.Range("E16").Value = '  -7.76%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001019'
$ws.Range("E17").Value = '  -6.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05709'
$ws.Range("E18").Value = '  -13.81%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9981'
$ws.Range("E19").Value = '  -0.44%  '

$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.18'
$ws.Range("E20").Value = '  -15.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.619'
$ws.Range("E21").Value = '  -8.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.76'
$ws.Range("E22").Value = '  -4.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.13'
$ws.Range("E23").Value = '  +3.85%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.264'
$ws.Range("E24").Value = '  -4.97%  '

$ws.Range("D25").Value = '20.072.62'
$ws.Range("E25").Value = '  -6.99%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.270'
$ws.Range("E26").Value = '  -3.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '135.29'
$ws.Range("E27").Value = '  -10.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.16'
$ws.Range("E28").Value = '  -6.78%  '

$ws.Range("D29").Value = '1.578.61'
$ws.Range("E29").Value = '  -7.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '109.88'
$ws.Range("E30").Value = '  -5.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.964'
$ws.Range("E31").Value = '  -18.90%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.371'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8305'
$ws.Range("E33").Value = '  -12.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07705'
$ws.Range("E34").Value = '  -3.81%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.494'
$ws.Range("E35").Value = '  +0.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.466'
$ws.Range("E36").Value = '  -1.44%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05901'
$ws.Range("E37").Value = '  -0.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.890'
$ws.Range("E38").Value = '  -5.26%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9973'
$ws.Range("E39").Value = '  -0.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '10.66'
$ws.Range("E40").Value = '  -4.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.02072'
$ws.Range("E41").Value = '  -6.01%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1915'
$ws.Range("E42").Value = '  -5.59%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.095'
$ws.Range("E43").Value = '  -6.95%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5312'
$ws.Range("E44").Value = '  -7.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.29'
$ws.Range("E45").Value = '  -5.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.531'
$ws.Range("E46").Value = '  -5.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5165'
$ws.Range("E47").Value = '  -6.75%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '114.56'
$ws.Range("E48").Value = '  -1.03%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.770'
$ws.Range("E49").Value = '  -6.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.044'
$ws.Range("E50").Value = '  -9.60%  '

$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9969'
$ws.Range("E51").Value = '  -0.58%  '
